# Update cryptocurrency price/volume data on Sheet1 to reflect the latest
# GitHub Actions scrape, including two rows whose coins swapped rank order.
#
# Price cells (column D) hold text that looks numeric (e.g. "245.26",
# "1.00", "2.006.50"); Excel would otherwise silently coerce such strings
# into numbers and drop formatting (trailing zeros, thousands separators).
# Setting NumberFormat to "@" (Text) first keeps them as text, matching the
# original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($row, $value) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

function Set-Volume($row, $value) {
    $ws.Cells.Item($row, 5).Value = $value
}

Set-Price 2 "36.874.55"
Set-Volume 2 "  -0.60%  "

Set-Price 3 "2.045.75"
Set-Volume 3 "  -0.32%  "

Set-Volume 4 "  +0.04%  "

Set-Price 5 "245.26"
Set-Volume 5 "  -1.88%  "

Set-Price 6 "0.654"
Set-Volume 6 "  -2.21%  "

Set-Price 7 "57.52"
Set-Volume 7 "  -3.72%  "

Set-Volume 8 "  +0.01%  "

Set-Price 9 "0.369"
Set-Volume 9 "  -4.93%  "

Set-Price 10 "0.0774"
Set-Volume 10 "  -2.55%  "

Set-Volume 11 "  +1.52%  "

Set-Price 12 "15.18"
Set-Volume 12 "  -5.74%  "

Set-Price 13 "0.869"
Set-Volume 13 "  +4.10%  "

Set-Price 14 "2.348.18"
Set-Volume 14 "  -0.12%  "

# --- Rows 15 & 16: Polkadot and WrappedEther swapped positions -------------

$ws.Cells.Item(15, 2).Value = "Polkadot"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-Price 15 "5.60"
Set-Volume 15 "  -3.46%  "

$ws.Cells.Item(16, 2).Value = "WrappedEther"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-Price 16 "2.006.50"
Set-Volume 16 "  -2.25%  "

# --- Continue simple updates -------------------------------------------

Set-Price 17 "17.88"
Set-Volume 17 "  -2.58%  "

Set-Price 18 "36.819.47"
Set-Volume 18 "  -0.70%  "

Set-Price 19 "73.28"
Set-Volume 19 "  -3.64%  "

Set-Price 20 "0.0₃0884"
Set-Volume 20 "  -2.48%  "

Set-Price 21 "5.37"
Set-Volume 21 "  -0.78%  "

Set-Price 22 "235.43"
Set-Volume 22 "  -1.30%  "

Set-Price 24 "2.44"
Set-Volume 24 "  +0.80%  "

Set-Price 25 "10.24"
Set-Volume 25 "  +8.69%  "

Set-Price 26 "2.18"
Set-Volume 26 "  -1.72%  "

Set-Price 27 "168.06"
Set-Volume 27 "  -0.95%  "

Set-Price 28 "19.91"
Set-Volume 28 "  -1.52%  "

Set-Volume 29 "  +13.15%  "

Set-Volume 30 "  -2.55%  "

Set-Price 31 "1.09"
Set-Volume 31 "  -4.64%  "

Set-Price 32 "4.67"
Set-Volume 32 "  +2.20%  "

Set-Price 33 "0.0609"
Set-Volume 33 "  -3.95%  "

# --- Rows 34 & 35: LidoDAOToken and BinanceUSD swapped positions -----------

$ws.Cells.Item(34, 2).Value = "BinanceUSD"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-Price 34 "1.00"
Set-Volume 34 "  +0.06%  "

$ws.Cells.Item(35, 2).Value = "LidoDAOToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-Price 35 "2.32"
Set-Volume 35 "  +3.32%  "

# --- Continue simple updates -------------------------------------------

Set-Volume 36 "  +4.23%  "

Set-Price 37 "0.0820"
Set-Volume 37 "  -8.11%  "

Set-Volume 38 "  -3.12%  "

Set-Price 39 "5.12"
Set-Volume 39 "  -3.26%  "

Set-Volume 40 "  -5.22%  "

Set-Price 41 "0.0222"
Set-Volume 41 "  -1.31%  "

Set-Volume 42 "  -0.42%  "

Set-Price 43 "0.0941"
Set-Volume 43 "  -13.99%  "

Set-Price 44 "96.30"
Set-Volume 44 "  -1.17%  "

Set-Price 45 "16.79"
Set-Volume 45 "  -4.67%  "

Set-Price 46 "1.294.56"
Set-Volume 46 "  -0.12%  "

Set-Price 47 "2.33"
Set-Volume 47 "  -7.45%  "

Set-Price 48 "2.85"
Set-Volume 48 "  -1.13%  "

Set-Price 49 "6.72"
Set-Volume 49 "  -2.28%  "

Set-Price 50 "2.232.92"
Set-Volume 50 "  -0.41%  "

Set-Price 51 "44.39"
Set-Volume 51 "  +0.17%  "
